# Add a new "Driver License" mapping entry to the IEPD mapping sheet.
#
# The existing sheet uses blank "spacer" rows between each logical
# block of mapping rows (row 40 is blank, row 41 starts the next
# block, etc). This change inserts a brand-new two-row block
# (one data row + one trailing blank spacer row) right before the
# "Registrant Residence Location" block, which pushes that block and
# everything below it down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank rows at row 41 (shifts old row 41.. down to 43..)
$ws.Rows("41:42").Insert()

# Populate the new row 41 with the Driver License mapping.
# Column C is written before column B so the shared-string table ends
# up with the same new-string ordering as the authored workbook
# (nc:IdentificationType before "Driver License ID").
$ws.Range("A41").Value = "Driver License"
$ws.Range("C41").Value = "nc:IdentificationType"
$ws.Range("B41").Value = "Driver License ID"
$ws.Range("D41").Value = "nc:IdentificationID"
$ws.Range("E41").Value = "niem-xsd:string"
$ws.Range("F41").Value = "exchange:FirearmRegistrationQueryResults/nc:DriverLicense/nc:DriverLicenseIdentification/nc:IdentificationID"
